# Individual_Reflection.docx edit
#  - justify (both) the existing body paragraphs
#  - tweak wording in paragraphs 2 and 3
#  - append two new body paragraphs after paragraph 3

$d = $word.ActiveDocument

function Replace-Text($scope, $findText, $replaceText) {
    $ok = $scope.Find.Execute($findText, $true, $false, $false, $false, $false,
                               $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# ---------------------------------------------------------------
# 1. Paragraph 2 ("In this project exercise ... optimal distance
#    wise.") - justify + insert a comma before " distance wise."
# ---------------------------------------------------------------
$d.Paragraphs(2).Format.Alignment = 3

Replace-Text $d.Content "that were optimal distance wise." "that were optimal, distance wise."

# ---------------------------------------------------------------
# 2. Paragraph 3 ("Again, the algorithms ability ...") - justify +
#    add possessive apostrophe + drop the trailing space.
# ---------------------------------------------------------------
$d.Paragraphs(3).Format.Alignment = 3

Replace-Text $d.Content "Again, the algorithms ability to calculate" "Again, the algorithms’ ability to calculate"
Replace-Text $d.Content "other greedy algorithms like Dijkstra’s. " "other greedy algorithms like Dijkstra’s."

# ---------------------------------------------------------------
# 3. Two brand-new paragraphs inserted right after paragraph 3,
#    before the trailing blank paragraph. Both inherit the same
#    pPr/rPr (Times New Roman, 24 half-points, double spacing) via
#    InsertParagraphAfter(), then get explicit justification.
# ---------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()
$pA = $d.Paragraphs(4)
$pA.Format.Alignment = 3
$pA.Range.Text = " My program’s implementation of the algorithm makes use of object instances of a custom Airport class. These are further represented as graph nodes and enqueued unto the search frontier which is ideally a priority queue (min heap) to take advantage of its sorted nature. A link or route between airports is represented as a route object, which is used to generate successor airports(nodes) in the search. These nodes on the frontier are arranged in the frontier based on a their distances to the destination airport. This enables airports with shorter distances to the destination airport to be expanded at the expense of others, hence enforcing that the shortest path is always returned by the algorithm. "

$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs(5)
$pB.Format.Alignment = 3
$pB.Range.Text = "Through this activity, I became well versed with the variations of searching algorithms and how they suite different needs. Breadth first search and depth first search are great graph search algorithms when cost is irrelevant to the search under question. However, in instances where optimality by a cost factor is required, greedy algorithms like Dijkstra’s and best first approaches are undisputably better options."

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
